$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The underlying shared-strings table got reshuffled and a new account
# ("000834301" / MARCUS) was inserted into the balances list, which
# cascades into rows 6-22 (everything else in the sheet is unchanged).
# Write out the resulting CONTA / NOME / SALDO values directly.

$data = @(
    @(6,  "004630773", "NABOR",       30667.24),
    @(7,  "004641487", "LAILA",       22063.42),
    @(8,  "004574428", "GUILHERME",   20365.37),
    @(9,  "004398253", "EULER",       19328.42),
    @(10, "004643737", "LARA",        17492.96),
    @(11, "000834301", "MARCUS",      15556.79),
    @(12, "004454365", "RAFAEL",      13566.65),
    @(13, "005064129", "THIAGO",      11241.46),
    @(14, "004927044", "CINTIA",      11053.46),
    @(15, "005020439", "BEATRIZ",     7369.3),
    @(16, "004346716", "TIAGO",       6395.86),
    @(17, "004206790", "EMMANUELLE",  5340.96),
    @(18, "000989247", "ANA",         4719.51),
    @(19, "004752519", "MARCUS",      4632.05),
    @(20, "004460491", "PEDRO",       3744.76),
    @(21, "004279859", "ASSOCIACAO",  2691.13),
    @(22, "004999410", "SONIA",       2420.16)
)

foreach ($row in $data) {
    $r = $row[0]
    $cellA = $ws.Cells.Item($r, 1)
    # Leading zeros must be preserved as text; a leading apostrophe forces
    # Excel to treat the numeric-looking string as text (like the original
    # file), then reset the style so no quote-prefix formatting sticks.
    $cellA.Value = "'" + $row[1]
    $cellA.Style = "Normal"
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
}

$ws.Range("N20").Select()
